$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the "Objetivos:" content in row 10 (was holding misplaced Teresa Paiva text) ---
$ws.Range("B10").Value = "Introduzir conceitos teóricos e práticos de Ecotoxicologia Aquática para estudantes de Engenharia Ambiental."
$ws.Range("C10").Value = "Introduzir conceitos teóricos e práticos de Ecotoxicologia Aquática para estudantes de Engenharia Ambiental."

# --- 2. Insert a new row at 13 to hold the "Docentes responsáveis:" value (Teresa Paiva), ---
#        pushing the rest of the table down by one row. Excel's native row-insert shifts
#        row heights/formats along with the content, so rows 13-21 retain the correct
#        per-row heights once shifted to 14-22.
$ws.Rows.Item(13).Insert()

# New row 13: B13/C13 hold the responsible-professor value, no label in A13, default row height.
# Copy number-format/font/alignment from row 14 (B/C) so the new cells use the existing
# column styles (s="2"/s="3") instead of picking up a stray default.
$ws.Range("B14:C14").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$ws.Range("B13").Value = "1720367 - Teresa Cristina Brazil de Paiva"
$ws.Range("C13").Value = "1720367 - Teresa Cristina Brazil de Paiva"
$ws.Range("A13").Clear()

# --- 3. Row 14 (was 13): "Programa resumido:" - replace placeholder "Semestral" text ---
$ws.Range("B14").Value = "Dar conhecimentos aos alunos de noções básicas teóricas e práticas sobre ecotoxicologia aquática e das técnicas usadas em laboratório para os cultivos e os ensaios com os organismos-teste padronizados."
$ws.Range("C14").Value = "Dar conhecimentos aos alunos de noções básicas teóricas e práticas sobre ecotoxicologia aquática e das técnicas usadas em laboratório para os cultivos e os ensaios com os organismos-teste padronizados."

# --- 4. Row 16 (was 15): "Programa:" - replace placeholder "01/01/2020" text ---
$ws.Range("B16").Value = "A. Teórico: 1. Ecotoxicologia: Introdução, histórico, conceitos; 2. Introdução de agentes químicos no ambiente aquático: biodisponibilidade de contaminantes, efeitos sinérgicos e antagônicos, impactos sobre os sistemas aquáticos; 3.Métodos de ensaios de toxicidade com organismos aquáticos: uso de bioindicadores; B. Prática: 4.Seleção, manutenção e cultivo de organismos aquáticos: boas práticas; 5. 5. Testes de toxicidade com organismos aquáticos."
$ws.Range("C16").Value = "A. Teórico: 1. Ecotoxicologia: Introdução, histórico, conceitos; 2. Introdução de agentes químicos no ambiente aquático: biodisponibilidade de contaminantes, efeitos sinérgicos e antagônicos, impactos sobre os sistemas aquáticos; 3.Métodos de ensaios de toxicidade com organismos aquáticos: uso de bioindicadores; B. Prática: 4.Seleção, manutenção e cultivo de organismos aquáticos: boas práticas; 5. 5. Testes de toxicidade com organismos aquáticos."

# --- 5. Row 19 (was 18): "Método:" - replace misplaced Teresa Paiva text ---
$ws.Range("B19").Value = "Aulas teóricas e práticas. Avaliação baseada em prova, exercício e relatório."
$ws.Range("C19").Value = "Aulas teóricas e práticas. Avaliação baseada em prova, exercício e relatório."

# --- 6. Row 20 (was 19): "Critério:" - update text (shifted up from old Norma line) ---
$ws.Range("B20").Value = "Média ponderada das notas atribuídas à prova, exercício e relatório."
$ws.Range("C20").Value = "Média ponderada das notas atribuídas à prova, exercício e relatório."

# --- 7. Row 21 (was 20): "Norma de recuperação:" - update text (shifted up from old Bibliografia line) ---
$ws.Range("B21").Value = "Nota final: NF ≥ 5,0"
$ws.Range("C21").Value = "Nota final: NF ≥ 5,0"

# --- 8. Row 22 (new): "Bibliografia:" with full reference text. Copy formats from row 21 ---
#        so the new row re-uses the existing column styles (s="1"/"2"/"3") instead of
#        creating fresh ones, then overwrite values and row height.
$ws.Range("A21:C21").Copy()
$ws.Range("A22:C22").PasteSpecial(-4122)
$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B22").Value = "Zagatto, P. A.; Bertoletti, E. Ecotoxicologia aquática– princípios e aplicações. RiMa, 2008.Azevedo, F.A.; Chasin, A.M. As bases toxicológicas da ecotoxicologia. RiMa, 2003.MOZETO, A. A.; UMBUZEIRO, G. A.; JARDIM, W. F. Métodos de coleta, análises físico-químicas e ensaios biológicos e ecotoxicológicos de sedimentos de água doce. São Carlos – SP. Cubo Multimídia & Propaganda, 2006."
$ws.Range("C22").Value = "Zagatto, P. A.; Bertoletti, E. Ecotoxicologia aquática– princípios e aplicações. RiMa, 2008.Azevedo, F.A.; Chasin, A.M. As bases toxicológicas da ecotoxicologia. RiMa, 2003.MOZETO, A. A.; UMBUZEIRO, G. A.; JARDIM, W. F. Métodos de coleta, análises físico-químicas e ensaios biológicos e ecotoxicológicos de sedimentos de água doce. São Carlos – SP. Cubo Multimídia & Propaganda, 2006."
$ws.Rows.Item(22).RowHeight = 120

